# Sujeto_4/Carbohidrates.xlsx -- "Add files via upload"
#
# The uploaded version replaces a handful of "Value (g)" column cells that
# had been entered as raw numbers (some of them scaled by 1000 and formatted
# with a thousands-separator number format) with plain decimal text, matching
# how every other row in the column is already stored (shared-string text).
#
# Rows (B column) changed from numeric to text:
#   B5  : 29385  (#,##0 style) -> "29.39"
#   B9  : 41               -> "41.0"
#   B10 : 15               -> "15.0"
#   B12 : 56               -> "56.0"
#   B17 : 16175  (#,##0 style) -> "16.18"
#   B19 : 37               -> "37.0"
#   B21 : 0                -> "0.0"
#   B22 : 15               -> "15.0"
#   B28 : 28175  (#,##0 style) -> "28.17"
#   B30 : 24               -> "24.0"
#   B31 : 120               -> "120.0"
#
# Setting .Value to a string that *looks* like a number still gets parsed
# back into a number by the Excel object model, so each cell is first forced
# into Text number-format, given its value, and then has that format cleared
# again (ClearFormats) so the stored cell ends up with the default style --
# just like every other text cell in the column -- while keeping the value
# stored as text (t="s").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "B5"  "29.39"
Set-TextValue "B9"  "41.0"
Set-TextValue "B10" "15.0"
Set-TextValue "B12" "56.0"
Set-TextValue "B17" "16.18"
Set-TextValue "B19" "37.0"
Set-TextValue "B21" "0.0"
Set-TextValue "B22" "15.0"
Set-TextValue "B28" "28.17"
Set-TextValue "B30" "24.0"
Set-TextValue "B31" "120.0"

# Reset the saved cursor/selection back to the top-left cell (the refreshed
# workbook no longer pins the view at D7).
$ws.Range("A1").Select() | Out-Null
